$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4: title/link update
$ws.Range("D4").Value = "PyTorch 이미지 데이터셋(Image Dataset) 구성에 관한 거의 모든 것!"
$ws.Range("E4").Value = "https://teddylee777.github.io/pytorch/dataset-dataloader"

# Row 51: title/link update
$ws.Range("D51").Value = "상업적 용도로 사용해도 괜찮은 CC0(Creative Commons 0) 라이센스"
$ws.Range("E51").Value = "https://bskyvision.com/entry/%EC%83%81%EC%97%85%EC%A0%81-%EC%9A%A9%EB%8F%84%EB%A1%9C-%EC%82%AC%EC%9A%A9%ED%95%B4%EB%8F%84-%EA%B4%9C%EC%B0%AE%EC%9D%80-CC0Creative-Commons-0-%EB%9D%BC%EC%9D%B4%EC%84%BC%EC%8A%A4"

# Row 52: title/link update
$ws.Range("D52").Value = "Relative Risk Regression"
$ws.Range("E52").Value = "http://ds.sumeun.org/?p=2603"
